# Apply the Team_PER_1995 update:
#  - Column B (Team) entries get reshuffled to a new team order
#  - Column C (PER) values get replaced with newly computed per-minute stats
# (per commit message: "fixed PER bug and added python ML code")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$teams = @{
    2  = "POR"
    3  = "NJN"
    4  = "CLE"
    5  = "DAL"
    6  = "MIA"
    7  = "SEA"
    8  = "ATL"
    9  = "MIL"
    10 = "LAC"
    11 = "DET"
    12 = "SAS"
    13 = "ORL"
    14 = "UTA"
    15 = "HOU"
    16 = "DEN"
    17 = "LAL"
    18 = "GSW"
    19 = "IND"
    20 = "CHI"
    21 = "PHI"
    22 = "CHH"
    23 = "BOS"
    24 = "WSB"
    25 = "SAC"
    26 = "PHO"
    27 = "NYK"
    28 = "MIN"
}

$values = @{
    2  = 13.35833333333333
    3  = 9.693333333333332
    4  = 14.4
    5  = 13.475
    6  = 13.71818181818182
    7  = 12.39285714285714
    8  = 10.75714285714286
    9  = 12.15
    10 = 12.29375
    11 = 13.13571428571429
    12 = 12.16428571428571
    13 = 15.02666666666667
    14 = 14.61666666666667
    15 = 12.26153846153846
    16 = 13.94117647058824
    17 = 13.94666666666667
    18 = 13.5625
    19 = 13.55714285714286
    20 = 13.88666666666666
    21 = 11.26875
    22 = 12.36
    23 = 13.72142857142857
    24 = 12.34
    25 = 11.26428571428572
    26 = 14.60666666666666
    27 = 11.36428571428571
    28 = 10.78571428571428
}

for ($row = 2; $row -le 28; $row++) {
    $ws.Cells.Item($row, 2).Value = $teams[$row]
    $ws.Cells.Item($row, 3).Value = $values[$row]
}
